$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the hidden "_GoBack" bookmark that sits right after the
#    "Sjelden problemer..." paragraph's text. The Bookmarks collection
#    does not surface this hidden bookmark (names starting with "_"
#    are hidden), so instead we locate the paragraph by its text,
#    delete that text run and reinsert it -- this drops the
#    bookmarkStart/bookmarkEnd anchors that lived at the end of the
#    deleted range, without touching anything else.
# ------------------------------------------------------------------
$goBackAnchor = $d.Content
$goBackAnchor.Find.Execute(
    "Sjelden problemer med denne*uttrekk fra en vgs).",
    $false, $false, $true, $false, $false, $true, 1, $false, "", 0
) | Out-Null

$goBackText = $goBackAnchor.Text
$goBackAnchor.Delete()
$goBackAnchor.InsertBefore($goBackText)

# ------------------------------------------------------------------
# 2) Insert a new empty paragraph right before the table, i.e. right
#    after the "Journalenheter beskriver..." paragraph.
# ------------------------------------------------------------------
$beforeTableAnchor = $d.Content
$beforeTableAnchor.Find.Execute(
    "Journalenheter beskriver hvem*papirbaserte arkiver.",
    $false, $false, $true, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$beforeTableAnchor.Collapse(0)  # wdCollapseEnd
$beforeTableAnchor.InsertParagraphAfter()

# ------------------------------------------------------------------
# 3) Append new paragraphs after the table (after the existing
#    trailing empty paragraph), at the very end of the document body:
#       AND/OR
#       (empty)
#       Varsel: Over 90% av journalenhetene er av samme type/navn.
#       (empty)
#       AND/OR
#       (empty)
#       Ingen journalenheter funnet.
#    NOTE: deliberately avoid touching Tables.Item(1).Range here --
#    reading a Table's Range taints subsequent Paragraphs lookups in
#    this host, so all positions are derived fresh from
#    Document.Content.End instead.
# ------------------------------------------------------------------
$lines = @(
    "AND/OR",
    "",
    "Varsel: Over 90% av journalenhetene er av samme type/navn.",
    "",
    "AND/OR",
    "",
    "Ingen journalenheter funnet."
)

foreach ($line in $lines) {
    $pos = $d.Content.End
    $r = $d.Range($pos, $pos)
    $r.InsertParagraphAfter()
    if ($line -ne "") {
        $pos2 = $d.Content.End - 1
        $r2 = $d.Range($pos2, $pos2)
        $r2.InsertBefore($line)
    }
}

Write-Output "done"
